$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two rows (25 and 26) first, so remaining rows keep their positions
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(25).Delete()

# New "Joueurs" text blocks, one per row (rows 2-24 correspond to B2:B24)
$players = @(
    "Saucisson, Ranzyo_S5, Thib_S5, Grenzo, Jilink, Chopa, Cha, Alan_S9",
    "The soulless, Killian, Julien_S7, Waikato, Mickey, Pilouche",
    "Chreet, Joshua, Gririsu, Hugo, Thib, Julien, Mickey_S8, Sunka",
    "Alecxperdu, Nelson, DTP, DramaPanda, Maitredoudou, FooD_Flo, Chreet_S5",
    "Nikolas, Gwendal, Clem fair play, Timeo, Amaury, Malkovan",
    "Bilel_S4, Bylost_S4, And, Xiao, Clovis, Thynael, Cosmos, Z4tix_S9",
    "Nelson_S4, Alex., Skowa_S4, Raphi, Xiao_S5, Manta, Dahmi1 Arti, Gobou_S6, Gwendal_S7, Lilian, Mehdiiii, Corentin",
    "Black Lolo, FooD_Flo_S5, Enze, Doggydog_S6, Enze_S7, Tim, Luc",
    "Kamiga, Z4tix, Mark, Alecxperdu_S2, Gobou, Polo, Skowa, Mark_S2, Schweppes",
    "Bylost, Nikolas_S5, Gririsu_S5, Nath_S6, Piiskoo",
    "Schweppes_S6, Ethan_S6, Gigi, Quentin, Enze_S8, Alex",
    "Eneko, Xori, Thyx",
    "Ethan, Nicovid, Ladoly, Didine, Deku, Sayo, Line_S5, Clem fair play_S5, TimeoGnc",
    "Cyrf, Sneus, Malkovan_S3, Julien_S5, Timeo_S5, Black Lolo_S5, Amaury_S5, Kwinn",
    "Line, HiYoucef, Saucisson_S5, Malkovan_S5, Vah Balress",
    "Crypto, Yatho, Mistimat, Spider, Guigui_S9",
    "Yuu, Nounours, Saminette, Nicovid_S6, Cyrf_S6, Blacks Star, Kamiga_S9",
    "Krak, Nwog, Jehovah, Alexadventure, Step, Feiik, Tiff, Natoxe, Bilal, Flau_S4, Bilel_S5",
    "Flau, Nath, Ethan_S3, Ranzyo, Guigui, Angel, Mtking, Dragon",
    "Bilel, Joshua_S3, GiulfeuYT, Doggydog, Alan, Mark_S5, Maitredoudou_S5",
    "Kamiga_S5, GiulfeuYT_S6, Alan_S7, Sy_boulette, Reiko, Ju, Mielle",
    "Blgham",
    "Lili, Samson, SorciShoot, Louan, Xo, Napoleon, Cyrf_S4, Espoir Perdu, And_S5, Major Chris, Nikolas_S6, Nemocca, Aurel, Theo, Chatoon, Clement, Mizuki, Armand, Manu"
)

# New "ID_Neurone" numbers for rows 2-24
$ids = @(1,2,3,4,5,6,7,8,10,11,12,13,15,16,17,18,19,20,21,22,23,24,25)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = $players[$i]
}
